$d = $word.ActiveDocument

# Title (appears twice: main heading + bold line near the end) - global replace
$d.Content.Find.Execute(
    "Play Cashpot Kegs Free: Unique Slot Game with Exciting Features", $true, $false, $false, $false, $false,
    $true, 1, $false, "Play Cashpot Kegs - Free Slot Game Review", 2
)

# "What we like" bullet list
$d.Content.Find.Execute(
    "Engaging gameplay and simple design", $true, $false, $false, $false, $false,
    $true, 1, $false, "Unique gameplay with classic grid structure", 2
)

$d.Content.Find.Execute(
    "Exciting bonus features with free spins and multipliers", $true, $false, $false, $false, $false,
    $true, 1, $false, "Interesting betting options and multipliers", 2
)

$d.Content.Find.Execute(
    "Cashpot feature adds to the excitement and chance to win big", $true, $false, $false, $false, $false,
    $true, 1, $false, "Exciting bonus game feature with free spins and wild multipliers", 2
)

$d.Content.Find.Execute(
    "Different betting options with multipliers can increase jackpot value", $true, $false, $false, $false, $false,
    $true, 1, $false, "Cashpot feature adds excitement and potential for big wins", 2
)

# "What we don't like" bullet list
$d.Content.Find.Execute(
    "Below average RTP of 95%", $true, $false, $false, $false, $false,
    $true, 1, $false, "RTP is slightly below average at 95%", 2
)

$d.Content.Find.Execute(
    "Minimal graphics and basic symbols", $true, $false, $false, $false, $false,
    $true, 1, $false, "Minimal graphics with basic symbols", 2
)

# Meta description (italic line)
$d.Content.Find.Execute(
    "Read our review of Cashpot Kegs, a unique slot game with exciting bonus features and a Cashpot jackpot. Play for free and discover your chance to win big.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Play Cashpot Kegs for free and discover unique gameplay, exciting features, and big win potential.", 2
)
